$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.744.71'
$ws.Range("E2").Value = '  -2.71%  '
$ws.Range("D3").Value = '2.399.87'
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.19%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.33%  '
$ws.Range("D9").Value = '2.380.45'
$ws.Range("E9").Value = '  -2.85%  '
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("E12").Value = '  -3.38%  '
$ws.Range("E13").Value = '  -2.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.77%  '
$ws.Range("E15").Value = '  -2.57%  '
$ws.Range("D17").Value = '60.668.54'
$ws.Range("E17").Value = '  -2.50%  '
$ws.Range("D18").Value = '2.391.94'
$ws.Range("E18").Value = '  -2.48%  '
$ws.Range("E19").Value = '  -3.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.46%  '
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '576.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.53%  '
$ws.Range("D30").Value = '0.0₃0905'
$ws.Range("E30").Value = '  -6.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.75'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.68%  '
$ws.Range("E32").Value = '  -7.23%  '
$ws.Range("E33").Value = '  -3.27%  '
$ws.Range("E34").Value = '  -8.06%  '
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.05%  '
$ws.Range("E37").Value = '  -3.62%  '
$ws.Range("E38").Value = '  -4.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '147.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("E41").Value = '  -5.25%  '
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("E44").Value = '  -5.39%  '
$ws.Range("E45").Value = '  -6.01%  '
$ws.Range("E46").Value = '  +20.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.90%  '
$ws.Range("E48").Value = '  -4.65%  '
$ws.Range("E49").Value = '  -3.43%  '
$ws.Range("E50").Value = '  -4.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.15%  '

Write-Output "Applied cryptos update"
